$d = $word.ActiveDocument

# 1) Heading: "Scalable Scalable Vector Graphics: ..." -> "Random-Access Rendering of General Vector Graphics"
$d.Content.Find.Execute(
    "Scalable Scalable Vector Graphics: Automatic Translation of Interactive SVGs to a Multithread VDOM for Fast Rendering",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Random-Access Rendering of General Vector Graphics", 2) | Out-Null

# 2) Paragraph: "This solution uses multithreaded CPU." -> new text (single run for now)
$d.Content.Find.Execute(
    "This solution uses multithreaded CPU.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Need to regenerate lattice mean shape must be static. Uses cell method (more details in  Precise Vector Textures for Real-Time 3D Rendering)", 2) | Out-Null

# 3) Delete the now-duplicate "Random-Access Rendering of General Vector Graphics" Heading2 paragraph
# (the one that originally sat right before "The A -buffer, an antialiased hidden surface method")
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Random-Access Rendering of General Vector Graphics`r") {
        $next = $para.Next()
        if ($next.Range.Text -eq "The A -buffer, an antialiased hidden surface method`r") {
            $para.Range.Delete()
            break
        }
    }
}

Write-Output "done"
